$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A, shifting existing columns (A:AC) right to (B:AD)
$ws.Columns("A").Insert()

# Header label for the new "Match ID" column (row 3 is the visible header row)
$ws.Range("A3").Value = "Match ID"

# Populate the new column with a constant Match ID of 1 for every data row (4-19)
# and for the totals row 20
$ws.Range("A4:A20").Value = 1

# Bold font to match the header styling used elsewhere on the sheet
$ws.Range("A3:A19").Font.Bold = $true

# Writing to the hidden totals row can leave a stray custom row height behind;
# auto-fit restores the sheet's default row height bookkeeping
$ws.Rows(20).AutoFit()

$ws.Range("A3:A20").Select()
